$d = $word.ActiveDocument

$pkgPrefix = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgSuffix = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Change 1: split "Convocatoria de Presentación: Junio 2022" into 3 runs, with
#     gramStart/gramEnd proofErr markers wrapping the standalone "Junio" run.
$frag1 = '<w:p><w:pPr><w:pStyle w:val="Default"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Convocatoria de Presentación: </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Junio</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> 2022</w:t></w:r></w:p>'
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Convocatoria de Presentación*") {
        $p.Range.InsertXML($pkgPrefix + $frag1 + $pkgSuffix)
        break
    }
}

# --- Change 2: give the "queríamos" paragraph 1.5-line spacing, drop the
#     spellStart/spellEnd proofErr wrapper, extend the sentence with new runs,
#     and leave a new empty paragraph behind it.
$frag2 = '<w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">Antes de empezar el proyecto no éramos realmente conscientes de la magnitud de un trabajo como </w:t></w:r><w:r><w:t xml:space="preserve">el que </w:t></w:r><w:r><w:t>queríamos por delante</w:t></w:r><w:r><w:t xml:space="preserve"> el cual era el desarrollo completo de una aplicación web, nuestras primeras expectaciones eran que iba a ser un trabajo relativamente sencillo poco a poco fueron tornándose a una visión más realista del proyecto </w:t></w:r></w:p><w:p/>'
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*queríamos*") {
        $p.Range.InsertXML($pkgPrefix + $frag2 + $pkgSuffix)
        break
    }
}
